$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The GO-term enrichment rows (2-7) were re-sorted into ascending adjusted
# p-value order (ties broken by GO id); row 1 (header) is unchanged.
# Each row keeps its original category/num_in_subset/num_total/adj_pval/
# term/ontology tuple - only the row order changes.

$ws.Range("A2").Value = "GO:0003735"
$ws.Range("B2").Value = 6
$ws.Range("C2").Value = 161
$ws.Range("D2").Value = 0.0498993422015866
$ws.Range("E2").Value = "structural constituent of ribosome"
$ws.Range("F2").Value = "MF"

$ws.Range("A3").Value = "GO:0005198"
$ws.Range("B3").Value = 8
$ws.Range("C3").Value = 190
$ws.Range("D3").Value = 0.0498993422015866
$ws.Range("E3").Value = "structural molecule activity"
$ws.Range("F3").Value = "MF"

$ws.Range("A4").Value = "GO:0005730"
$ws.Range("B4").Value = 13
$ws.Range("C4").Value = 307
$ws.Range("D4").Value = 0.0384649508592009
$ws.Range("E4").Value = "nucleolus"
$ws.Range("F4").Value = "CC"

$ws.Range("A5").Value = "GO:0009987"
$ws.Range("B5").Value = 179
$ws.Range("C5").Value = 2338
$ws.Range("D5").Value = 0.0498993422015866
$ws.Range("E5").Value = "cellular process"
$ws.Range("F5").Value = "BP"

$ws.Range("A6").Value = "GO:0043228"
$ws.Range("B6").Value = 66
$ws.Range("C6").Value = 1305
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = "non-membrane-bounded organelle"
$ws.Range("F6").Value = "CC"

$ws.Range("A7").Value = "GO:0043232"
$ws.Range("B7").Value = 66
$ws.Range("C7").Value = 1305
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = "intracellular non-membrane-bounded organelle"
$ws.Range("F7").Value = "CC"
